$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(58, 8).Value = 2471.5386  # H58: 1469.8462 -> 2471.5386
$ws.Cells.Item(58, 9).Value = 700  # I58: 369.77777 -> 700
$ws.Cells.Item(58, 10).Value = 3990  # J58: 3945 -> 3990
$ws.Cells.Item(58, 11).Value = 2100  # K58: 1109.33331 -> 2100
$ws.Cells.Item(58, 12).Value = 11970  # L58: 11835 -> 11970
$ws.Cells.Item(58, 13).Value = -1950  # M58: -959.33331 -> -1950
$ws.Cells.Item(58, 14).Value = -12270  # N58: -12135 -> -12270

$ws.Cells.Item(81, 8).Value = 40000  # H81: 0 -> 40000
$ws.Cells.Item(81, 10).Value = 40000  # J81: 0 -> 40000
$ws.Cells.Item(81, 12).Value = 40000  # L81: 0 -> 40000
$ws.Cells.Item(81, 14).Value = -41996  # N81: None -> -41996

$ws.Cells.Item(84, 8).Value = 40000  # H84: 0 -> 40000
$ws.Cells.Item(84, 10).Value = 40000  # J84: 0 -> 40000
$ws.Cells.Item(84, 12).Value = 120000  # L84: 0 -> 120000
$ws.Cells.Item(84, 14).Value = -129984  # N84: None -> -129984

$ws.Cells.Item(88, 8).Value = 418.16  # H88: 567.3333 -> 418.16
$ws.Cells.Item(88, 9).Value = 286.1111  # I88: 550 -> 286.1111
$ws.Cells.Item(88, 10).Value = 492.4375  # J88: 576 -> 492.4375
$ws.Cells.Item(88, 11).Value = 286.1111  # K88: 550 -> 286.1111
$ws.Cells.Item(88, 12).Value = 492.4375  # L88: 576 -> 492.4375
$ws.Cells.Item(88, 13).Value = 119.8889  # M88: -144 -> 119.8889
$ws.Cells.Item(88, 14).Value = -1304.4375  # N88: -1388 -> -1304.4375

$ws.Cells.Item(91, 8).Value = 418.16  # H91: 567.3333 -> 418.16
$ws.Cells.Item(91, 9).Value = 286.1111  # I91: 550 -> 286.1111
$ws.Cells.Item(91, 10).Value = 492.4375  # J91: 576 -> 492.4375
$ws.Cells.Item(91, 11).Value = 286.1111  # K91: 550 -> 286.1111
$ws.Cells.Item(91, 12).Value = 492.4375  # L91: 576 -> 492.4375
$ws.Cells.Item(91, 13).Value = 1117.8889  # M91: 854 -> 1117.8889
$ws.Cells.Item(91, 14).Value = -3300.4375  # N91: -3384 -> -3300.4375

$ws.Cells.Item(120, 8).Value = 36950  # H120: 37450 -> 36950
$ws.Cells.Item(120, 10).Value = 36950  # J120: 37450 -> 36950
$ws.Cells.Item(120, 12).Value = 36950  # L120: 37450 -> 36950
$ws.Cells.Item(120, 14).Value = -46626  # N120: -47126 -> -46626

$ws.Cells.Item(126, 8).Value = 52000  # H126: 46000 -> 52000
$ws.Cells.Item(126, 9).Value = 52000  # I126: 0 -> 52000
$ws.Cells.Item(126, 10).Value = 0  # J126: 46000 -> 0
$ws.Cells.Item(126, 11).Value = 52000  # K126: 0 -> 52000
$ws.Cells.Item(126, 12).ClearContents()  # L126: 46000 -> (removed)
$ws.Cells.Item(126, 13).Value = -47060  # M126: None -> -47060
$ws.Cells.Item(126, 14).Value = 0  # N126: -55880 -> 0

$ws.Cells.Item(135, 8).Value = 1341.8572  # H135: 1042.9259 -> 1341.8572
$ws.Cells.Item(135, 9).Value = 1402.1333  # I135: 1011.85 -> 1402.1333
$ws.Cells.Item(135, 10).Value = 1191.1666  # J135: 1131.7142 -> 1191.1666
$ws.Cells.Item(135, 11).Value = 12619.1997  # K135: 9106.65 -> 12619.1997
$ws.Cells.Item(135, 12).Value = 10720.4994  # L135: 10185.4278 -> 10720.4994
$ws.Cells.Item(135, 13).Value = -10084.1997  # M135: -6571.65 -> -10084.1997
$ws.Cells.Item(135, 14).Value = -15790.4994  # N135: -15255.4278 -> -15790.4994

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(11, 8).Value = 4000334.2  # H11: 3000300.8 -> 4000334.2
$ws.Cells.Item(11, 9).Value = 4000334.2  # I11: 3000300.8 -> 4000334.2
$ws.Cells.Item(11, 11).Value = 4000334.2  # K11: 3000300.8 -> 4000334.2
$ws.Cells.Item(11, 13).Value = -4000190.2  # M11: -3000156.8 -> -4000190.2

$ws.Cells.Item(61, 8).Value = 1963.8889  # H61: 2050.838 -> 1963.8889
$ws.Cells.Item(61, 9).Value = 2025.25  # I61: 2184.5264 -> 2025.25
$ws.Cells.Item(61, 10).Value = 1893.762  # J61: 1909.7222 -> 1893.762
$ws.Cells.Item(61, 11).Value = 2025.25  # K61: 2184.5264 -> 2025.25
$ws.Cells.Item(61, 12).Value = 1893.762  # L61: 1909.7222 -> 1893.762
$ws.Cells.Item(61, 13).Value = -1813.25  # M61: -1972.5264 -> -1813.25
$ws.Cells.Item(61, 14).Value = -2317.762  # N61: -2333.7222 -> -2317.762

$ws.Cells.Item(74, 8).Value = 3931.2  # H74: 1731.5714 -> 3931.2
$ws.Cells.Item(74, 9).Value = 5598.5713  # I74: 1595.9048 -> 5598.5713
$ws.Cells.Item(74, 10).Value = 1809.091  # J74: 1935.0714 -> 1809.091
$ws.Cells.Item(74, 11).Value = 5598.5713  # K74: 1595.9048 -> 5598.5713
$ws.Cells.Item(74, 12).Value = 1809.091  # L74: 1935.0714 -> 1809.091
$ws.Cells.Item(74, 13).Value = -4724.5713  # M74: -721.9048 -> -4724.5713
$ws.Cells.Item(74, 14).Value = -3557.091  # N74: -3683.0714 -> -3557.091

$ws.Cells.Item(77, 8).Value = 3931.2  # H77: 1731.5714 -> 3931.2
$ws.Cells.Item(77, 9).Value = 5598.5713  # I77: 1595.9048 -> 5598.5713
$ws.Cells.Item(77, 10).Value = 1809.091  # J77: 1935.0714 -> 1809.091
$ws.Cells.Item(77, 11).Value = 27992.8565  # K77: 7979.524 -> 27992.8565
$ws.Cells.Item(77, 12).Value = 9045.455  # L77: 9675.357 -> 9045.455
$ws.Cells.Item(77, 13).Value = -23624.8565  # M77: -3611.524 -> -23624.8565
$ws.Cells.Item(77, 14).Value = -17781.455  # N77: -18411.357 -> -17781.455

$ws.Cells.Item(102, 8).Value = 2330.0715  # H102: 3089 -> 2330.0715
$ws.Cells.Item(102, 9).Value = 1213  # I102: 1490.3334 -> 1213
$ws.Cells.Item(102, 10).Value = 3447.1428  # J102: 3888.3333 -> 3447.1428
$ws.Cells.Item(102, 11).Value = 1213  # K102: 1490.3334 -> 1213
$ws.Cells.Item(102, 12).Value = 3447.1428  # L102: 3888.3333 -> 3447.1428
$ws.Cells.Item(102, 13).Value = 409  # M102: 131.6666 -> 409
$ws.Cells.Item(102, 14).Value = -6691.1428  # N102: -7132.3333 -> -6691.1428

$ws.Cells.Item(136, 8).Value = 1963.8889  # H136: 2050.838 -> 1963.8889
$ws.Cells.Item(136, 9).Value = 2025.25  # I136: 2184.5264 -> 2025.25
$ws.Cells.Item(136, 10).Value = 1893.762  # J136: 1909.7222 -> 1893.762
$ws.Cells.Item(136, 11).Value = 6075.75  # K136: 6553.5792 -> 6075.75
$ws.Cells.Item(136, 12).Value = 5681.286  # L136: 5729.1666 -> 5681.286
$ws.Cells.Item(136, 13).Value = -3525.75  # M136: -4003.5792 -> -3525.75
$ws.Cells.Item(136, 14).Value = -10781.286  # N136: -10829.1666 -> -10781.286

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(8, 8).Value = 6833.3335  # H8: 10000 -> 6833.3335
$ws.Cells.Item(8, 10).Value = 2750  # J8: 5000 -> 2750
$ws.Cells.Item(8, 12).Value = 2750  # L8: 5000 -> 2750
$ws.Cells.Item(8, 14).Value = -3030  # N8: -5280 -> -3030

$ws.Cells.Item(86, 8).Value = 1216.7142  # H86: 1114.2858 -> 1216.7142
$ws.Cells.Item(86, 9).Value = 1022.125  # I86: 1200 -> 1022.125
$ws.Cells.Item(86, 10).Value = 1476.1666  # J86: 1000 -> 1476.1666
$ws.Cells.Item(86, 11).Value = 1022.125  # K86: 1200 -> 1022.125
$ws.Cells.Item(86, 12).Value = 1476.1666  # L86: 1000 -> 1476.1666
$ws.Cells.Item(86, 13).Value = 100.875  # M86: -77 -> 100.875
$ws.Cells.Item(86, 14).Value = -3722.1666  # N86: -3246 -> -3722.1666

$ws.Cells.Item(89, 8).Value = 1216.7142  # H89: 1114.2858 -> 1216.7142
$ws.Cells.Item(89, 9).Value = 1022.125  # I89: 1200 -> 1022.125
$ws.Cells.Item(89, 10).Value = 1476.1666  # J89: 1000 -> 1476.1666
$ws.Cells.Item(89, 11).Value = 5110.625  # K89: 6000 -> 5110.625
$ws.Cells.Item(89, 12).Value = 7380.833000000001  # L89: 5000 -> 7380.833000000001
$ws.Cells.Item(89, 13).Value = 505.375  # M89: -384 -> 505.375
$ws.Cells.Item(89, 14).Value = -18612.833  # N89: -16232 -> -18612.833

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(3, 8).Value = 33920924  # H3: 38160800 -> 33920924
$ws.Cells.Item(3, 9).Value = 1950  # I3: 2000 -> 1950
$ws.Cells.Item(3, 11).Value = 1950  # K3: 2000 -> 1950
$ws.Cells.Item(3, 13).Value = -1837  # M3: -1887 -> -1837

$ws.Cells.Item(43, 8).Value = 17890  # H43: 17862.637 -> 17890
$ws.Cells.Item(43, 10).Value = 17890  # J43: 17862.637 -> 17890
$ws.Cells.Item(43, 12).Value = 17890  # L43: 17862.637 -> 17890
$ws.Cells.Item(43, 14).Value = -18258  # N43: -18230.637 -> -18258

$ws.Cells.Item(101, 8).Value = 17890  # H101: 17862.637 -> 17890
$ws.Cells.Item(101, 10).Value = 17890  # J101: 17862.637 -> 17890
$ws.Cells.Item(101, 12).Value = 17890  # L101: 17862.637 -> 17890
$ws.Cells.Item(101, 14).Value = -24380  # N101: -24352.637 -> -24380

$ws.Cells.Item(138, 8).Value = 39463.332  # H138: 39464.168 -> 39463.332
$ws.Cells.Item(138, 10).Value = 41689.09  # J138: 41690 -> 41689.09
$ws.Cells.Item(138, 12).Value = 41689.09  # L138: 41690 -> 41689.09
$ws.Cells.Item(138, 14).Value = -51969.09  # N138: -51970 -> -51969.09

$ws.Cells.Item(139, 8).Value = 30283  # H139: 30949.666 -> 30283
$ws.Cells.Item(139, 10).Value = 33600  # J139: 34933.332 -> 33600
$ws.Cells.Item(139, 12).Value = 33600  # L139: 34933.332 -> 33600
$ws.Cells.Item(139, 14).Value = -43880  # N139: -45213.332 -> -43880

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, 8).Value = 185.6  # H23: 176.07143 -> 185.6
$ws.Cells.Item(23, 9).Value = 130.16667  # I23: 93 -> 130.16667
$ws.Cells.Item(23, 10).Value = 222.55556  # J23: 209.3 -> 222.55556
$ws.Cells.Item(23, 11).Value = 390.50001  # K23: 279 -> 390.50001
$ws.Cells.Item(23, 12).Value = 667.66668  # L23: 627.9000000000001 -> 667.66668
$ws.Cells.Item(23, 13).Value = -155.50001  # M23: -44 -> -155.50001
$ws.Cells.Item(23, 14).Value = -1137.66668  # N23: -1097.9 -> -1137.66668

$ws.Cells.Item(34, 8).Value = 1235.8572  # H34: 1590.8182 -> 1235.8572
$ws.Cells.Item(34, 10).Value = 1375.1666  # J34: 1855.4445 -> 1375.1666
$ws.Cells.Item(34, 12).Value = 4125.4998  # L34: 5566.333500000001 -> 4125.4998
$ws.Cells.Item(34, 14).Value = -4293.4998  # N34: -5734.333500000001 -> -4293.4998

$ws.Cells.Item(38, 8).Value = 231.36363  # H38: 194.875 -> 231.36363
$ws.Cells.Item(38, 9).Value = 292.14285  # I38: 260.625 -> 292.14285
$ws.Cells.Item(38, 10).Value = 125  # J38: 129.125 -> 125
$ws.Cells.Item(38, 11).Value = 876.4285500000001  # K38: 781.875 -> 876.4285500000001
$ws.Cells.Item(38, 12).Value = 375  # L38: 387.375 -> 375
$ws.Cells.Item(38, 13).Value = -529.4285500000001  # M38: -434.875 -> -529.4285500000001
$ws.Cells.Item(38, 14).Value = -1069  # N38: -1081.375 -> -1069

$ws.Cells.Item(48, 8).Value = 980  # H48: 0 -> 980
$ws.Cells.Item(48, 9).Value = 980  # I48: 0 -> 980
$ws.Cells.Item(48, 11).Value = 2940  # K48: 0 -> 2940
$ws.Cells.Item(48, 13).Value = -2690  # M48: None -> -2690

$ws.Cells.Item(51, 8).Value = 1248.4286  # H51: 1637.6666 -> 1248.4286
$ws.Cells.Item(51, 10).Value = 0  # J51: 3000 -> 0
$ws.Cells.Item(51, 12).Value = 0  # L51: 9000 -> 0
$ws.Cells.Item(51, 14).ClearContents()  # N51: -9920 -> (removed)

$ws.Cells.Item(57, 8).Value = 2970  # H57: 2813.3333 -> 2970
$ws.Cells.Item(57, 9).Value = 2970  # I57: 2813.3333 -> 2970
$ws.Cells.Item(57, 11).Value = 8910  # K57: 8439.999899999999 -> 8910
$ws.Cells.Item(57, 13).Value = -8351  # M57: -7880.999899999999 -> -8351

$ws.Cells.Item(93, 8).Value = 6398.5  # H93: 6400 -> 6398.5
$ws.Cells.Item(93, 10).Value = 6398.5  # J93: 6400 -> 6398.5
$ws.Cells.Item(93, 12).Value = 19195.5  # L93: 19200 -> 19195.5
$ws.Cells.Item(93, 14).Value = -22939.5  # N93: -22944 -> -22939.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3975.375  # H7: 4436.727 -> 3975.375
$ws.Cells.Item(7, 9).Value = 3650.375  # I7: 4374.75 -> 3650.375
$ws.Cells.Item(7, 10).Value = 4300.375  # J7: 4472.143 -> 4300.375
$ws.Cells.Item(7, 11).Value = 3650.375  # K7: 4374.75 -> 3650.375
$ws.Cells.Item(7, 12).Value = 4300.375  # L7: 4472.143 -> 4300.375
$ws.Cells.Item(7, 13).Value = -3538.375  # M7: -4262.75 -> -3538.375
$ws.Cells.Item(7, 14).Value = -4524.375  # N7: -4696.143 -> -4524.375

$ws.Cells.Item(40, 8).Value = 55559784  # H40: 58827916 -> 55559784
$ws.Cells.Item(40, 10).Value = 4694.0835  # J40: 4984.4546 -> 4694.0835
$ws.Cells.Item(40, 12).Value = 4694.0835  # L40: 4984.4546 -> 4694.0835
$ws.Cells.Item(40, 14).Value = -4966.0835  # N40: -5256.4546 -> -4966.0835

$ws.Cells.Item(46, 8).Value = 1667241.6  # H46: 1250556.2 -> 1667241.6

$ws.Cells.Item(68, 8).Value = 2100  # H68: 2118.4814 -> 2100
$ws.Cells.Item(68, 9).Value = 1700  # I68: 1674.9375 -> 1700
$ws.Cells.Item(68, 10).Value = 3500  # J68: 2763.6365 -> 3500
$ws.Cells.Item(68, 11).Value = 1700  # K68: 1674.9375 -> 1700
$ws.Cells.Item(68, 12).Value = 3500  # L68: 2763.6365 -> 3500
$ws.Cells.Item(68, 13).Value = -951  # M68: -925.9375 -> -951
$ws.Cells.Item(68, 14).Value = -4998  # N68: -4261.636500000001 -> -4998

$ws.Cells.Item(71, 8).Value = 2100  # H71: 2118.4814 -> 2100
$ws.Cells.Item(71, 9).Value = 1700  # I71: 1674.9375 -> 1700
$ws.Cells.Item(71, 10).Value = 3500  # J71: 2763.6365 -> 3500
$ws.Cells.Item(71, 11).Value = 8500  # K71: 8374.6875 -> 8500
$ws.Cells.Item(71, 12).Value = 17500  # L71: 13818.1825 -> 17500
$ws.Cells.Item(71, 13).Value = -4756  # M71: -4630.6875 -> -4756
$ws.Cells.Item(71, 14).Value = -24988  # N71: -21306.1825 -> -24988

$ws.Cells.Item(100, 8).Value = 3531.4119  # H100: 2411.7666 -> 3531.4119
$ws.Cells.Item(100, 9).Value = 4075.25  # I100: 1706.25 -> 4075.25
$ws.Cells.Item(100, 10).Value = 3364.077  # J100: 3218.0715 -> 3364.077
$ws.Cells.Item(100, 11).Value = 4075.25  # K100: 1706.25 -> 4075.25
$ws.Cells.Item(100, 12).Value = 3364.077  # L100: 3218.0715 -> 3364.077
$ws.Cells.Item(100, 13).Value = -3534.25  # M100: -1165.25 -> -3534.25
$ws.Cells.Item(100, 14).Value = -4446.077  # N100: -4300.0715 -> -4446.077

$ws.Cells.Item(126, 8).Value = 3975.375  # H126: 4436.727 -> 3975.375
$ws.Cells.Item(126, 9).Value = 3650.375  # I126: 4374.75 -> 3650.375
$ws.Cells.Item(126, 10).Value = 4300.375  # J126: 4472.143 -> 4300.375
$ws.Cells.Item(126, 11).Value = 10951.125  # K126: 13124.25 -> 10951.125
$ws.Cells.Item(126, 12).Value = 12901.125  # L126: 13416.429 -> 12901.125
$ws.Cells.Item(126, 13).Value = -8481.125  # M126: -10654.25 -> -8481.125
$ws.Cells.Item(126, 14).Value = -17841.125  # N126: -18356.429 -> -17841.125

$ws.Cells.Item(132, 8).Value = 31577.666  # H132: 44908.92 -> 31577.666
$ws.Cells.Item(132, 9).Value = 34703.062  # I132: 57227.367 -> 34703.062
$ws.Cells.Item(132, 10).Value = 6574.5  # J132: 5900.5 -> 6574.5
$ws.Cells.Item(132, 11).Value = 104109.186  # K132: 171682.101 -> 104109.186
$ws.Cells.Item(132, 12).Value = 19723.5  # L132: 17701.5 -> 19723.5
$ws.Cells.Item(132, 13).Value = -101579.186  # M132: -169152.101 -> -101579.186
$ws.Cells.Item(132, 14).Value = -24783.5  # N132: -22761.5 -> -24783.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 1488.3636  # H81: 1711.4348 -> 1488.3636
$ws.Cells.Item(81, 9).Value = 1466.5264  # I81: 1654.6111 -> 1466.5264
$ws.Cells.Item(81, 10).Value = 1626.6666  # J81: 1916 -> 1626.6666
$ws.Cells.Item(81, 11).Value = 2933.0528  # K81: 3309.2222 -> 2933.0528
$ws.Cells.Item(81, 12).Value = 3253.3332  # L81: 3832 -> 3253.3332
$ws.Cells.Item(81, 13).Value = -1872.0528  # M81: -2248.2222 -> -1872.0528
$ws.Cells.Item(81, 14).Value = -5375.3332  # N81: -5954 -> -5375.3332

$ws.Cells.Item(84, 8).Value = 1488.3636  # H84: 1711.4348 -> 1488.3636
$ws.Cells.Item(84, 9).Value = 1466.5264  # I84: 1654.6111 -> 1466.5264
$ws.Cells.Item(84, 10).Value = 1626.6666  # J84: 1916 -> 1626.6666
$ws.Cells.Item(84, 11).Value = 14665.264  # K84: 16546.111 -> 14665.264
$ws.Cells.Item(84, 12).Value = 16266.666  # L84: 19160 -> 16266.666
$ws.Cells.Item(84, 13).Value = -9361.263999999999  # M84: -11242.111 -> -9361.263999999999
$ws.Cells.Item(84, 14).Value = -26874.666  # N84: -29768 -> -26874.666

$ws.Cells.Item(122, 8).Value = 2829.6667  # H122: 1765.4667 -> 2829.6667
$ws.Cells.Item(122, 9).Value = 2752  # I122: 1470.7 -> 2752
$ws.Cells.Item(122, 10).Value = 2985  # J122: 2355 -> 2985
$ws.Cells.Item(122, 11).Value = 8256  # K122: 4412.1 -> 8256
$ws.Cells.Item(122, 12).Value = 8955  # L122: 7065 -> 8955
$ws.Cells.Item(122, 13).Value = -5806  # M122: -1962.1 -> -5806
$ws.Cells.Item(122, 14).Value = -13855  # N122: -11965 -> -13855

$ws.Cells.Item(132, 8).Value = 1678.4791  # H132: 1671.7693 -> 1678.4791
$ws.Cells.Item(132, 9).Value = 1399.0883  # I132: 1498.7142 -> 1399.0883
$ws.Cells.Item(132, 10).Value = 2357  # J132: 2028.0588 -> 2357
$ws.Cells.Item(132, 11).Value = 4197.2649  # K132: 4496.142599999999 -> 4197.2649
$ws.Cells.Item(132, 12).Value = 7071  # L132: 6084.1764 -> 7071
$ws.Cells.Item(132, 13).Value = -1667.2649  # M132: -1966.142599999999 -> -1667.2649
$ws.Cells.Item(132, 14).Value = -12131  # N132: -11144.1764 -> -12131

$ws.Cells.Item(133, 8).Value = 31000  # H133: 30000 -> 31000
$ws.Cells.Item(133, 10).Value = 31000  # J133: 30000 -> 31000
$ws.Cells.Item(133, 12).Value = 31000  # L133: 30000 -> 31000
$ws.Cells.Item(133, 14).Value = -41120  # N133: -40120 -> -41120

$ws.Cells.Item(136, 8).Value = 2271.5134  # H136: 2077.6562 -> 2271.5134
$ws.Cells.Item(136, 9).Value = 1779.0769  # I136: 1848.6957 -> 1779.0769
$ws.Cells.Item(136, 10).Value = 3435.4546  # J136: 2662.7778 -> 3435.4546
$ws.Cells.Item(136, 11).Value = 5337.2307  # K136: 5546.0871 -> 5337.2307
$ws.Cells.Item(136, 12).Value = 10306.3638  # L136: 7988.3334 -> 10306.3638
$ws.Cells.Item(136, 13).Value = -2787.2307  # M136: -2996.0871 -> -2787.2307
$ws.Cells.Item(136, 14).Value = -15406.3638  # N136: -13088.3334 -> -15406.3638
